$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2 (CasesTab query): append ORDER BY / LIMIT clause
$origB2 = $ws.Range("B2").Text
$ws.Range("B2").Value = $origB2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# B3 (SamplesTab query): append ORDER BY / LIMIT clause
$origB3 = $ws.Range("B3").Text
$ws.Range("B3").Value = $origB3 + "`n order By samp.sample_id ASC LIMIT 100"

# B4 (FilesTab query): replace trailing "order by" clause with ORDER BY / LIMIT clause
$origB4 = $ws.Range("B4").Text
$tailOld = "`n    order by f.file_name"
$tailNew = "`n     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value = $origB4.Replace($tailOld, $tailNew)

# Row heights grow by one wrapped line because the wrapped text is now longer
# (matches Excel's automatic row-height reflow after the text edit)
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Update the active selection/scroll position to match the author's final cursor position
$ws.Range("B4").Select() | Out-Null
